$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 13) for year 2021, mirroring the layout of the
# preceding rows: column A holds the year label (styled like the other year
# cells), columns B:O hold the numeric figures.

# Copy the formatting of the last existing data row (row 12) down into the
# new row so the new cells pick up the same style (border/alignment/bold).
$ws.Range("A12:O12").Copy()
$ws.Range("A13:O13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"

$values = @(1000, 5979, 228, 784, 403, 549, 366, 165, 58, 849, 417, 458, 539, 11795)
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 2 + $i  # Column B is index 2
    $ws.Cells.Item(13, $col).Value = $values[$i]
}
